$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder / relabel the header row.
# Write C1 ("Speg") first so the existing shared string stays referenced
# while it moves from D1 -> C1, then introduce the brand-new unique
# strings in the same order Excel would have recorded them.
$ws.Range("C1").Value = "Speg"
$ws.Range("A2").Value = "xyz"
$ws.Range("C3").Value = "spe"
$ws.Range("A1").Value = "Brand Name"
$ws.Range("B1").Value = "Alias"
$ws.Range("D1").Value = "Lpeg"
$ws.Range("E2").Value = "btl"
$ws.Range("A3").Value = "xabc"
$ws.Range("F1").Value = "License Name"

# Remaining values that don't introduce new unique strings.
$ws.Range("B2").Value = 1000
$ws.Range("F2").Value = "xyz"
$ws.Range("B3").Value = 750
$ws.Range("F3").Value = "xyz"

$ws.Range("F3").Select()
